$d = $word.ActiveDocument

# --- 1) Remove the trailing "," run that follows the second inline picture
#        in the paragraph whose visible text is ", ," (two formulas joined
#        by commas, rendered as two inline images around the text). ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -eq 2 -and $p.Range.Text -eq ", ,`r") {
        $r = $p.Range
        # Characters.Count includes the trailing paragraph mark, so the
        # last real character (the stray ",") is one before that.
        $last = $r.Characters.Item($r.Characters.Count - 1)
        if ($last.Text -eq ",") {
            $last.Delete()
        }
        break
    }
}

# --- 2) Move the lone "_GoBack" bookmark from the end of the
#        "k = 6 m=2.0 a = 0.5" paragraph to the empty paragraph that sits
#        right after "Изучить решение задачи Коши..." (i.e. right before
#        the "Задание" heading). Adding a new "_GoBack" bookmark implicitly
#        relocates the existing one, since it is a singleton bookmark. ---
for ($i = 2; $i -lt $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Length -eq 1) {
        $prev = $d.Paragraphs.Item($i - 1).Range.Text
        $next = $d.Paragraphs.Item($i + 1).Range.Text
        if ($prev -like "*Коши*" -and $next -like "*Задание*") {
            $d.Bookmarks.Add("_GoBack", $p.Range)
            break
        }
    }
}
